$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the 2020 column (L) with two new year columns: 2021 (M) and 2022 (N),
# matching the formatting already used by the existing year columns.
$ws.Range("L3:L4").Copy($ws.Range("M3")) | Out-Null
$ws.Range("L3:L4").Copy($ws.Range("N3")) | Out-Null

$ws.Range("M3").Value = 2021
$ws.Range("N3").Value = 2022
$ws.Range("M4").Value = 6.18
$ws.Range("N4").Value = 6.18

# Match the workbook's saved selection state.
$ws.Range("N15").Select() | Out-Null
